$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old trailing header cells (Branch / Address) that are no
# longer part of the table.
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""

# New headers
$ws.Range("A1").Value = "TC_ID"
$ws.Range("B1").Value = "TC_DESC"
$ws.Range("C1").Value = "TC_Status"

# Column A (TC_ID) values
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Column B (TC_DESC) values
$ws.Range("B2").Value = "qww"
$ws.Range("B3").Value = "eee"
$ws.Range("B4").Value = "fff"

# Column C (TC_Status) values
$ws.Range("C3").Value = "Yes"
$ws.Range("C4").Value = "No"
$ws.Range("C2").Value = "Passed"

# Restore the active selection to match the saved worksheet view.
$ws.Range("G8").Select() | Out-Null
